# Adjust MAMBI Salinity standards to deal with single TF stations
# Add two placeholder Tidal Freshwater rows (good/bad reference) to the
# "Saline Sites" standards sheet so the PCA step has multiple values to
# work with for lone tidal freshwater sites.

$wb = $excel.ActiveWorkbook

$salineSheet = $wb.Worksheets.Item("Saline Sites")

# --- Append the two placeholder rows to "Saline Sites" ---
# Row 16: TF_bad_placeholder
$salineSheet.Range("A16").Value = "TF_bad_placeholder"
$salineSheet.Range("B16").Value = 0
$salineSheet.Range("C16").Value = 6
$salineSheet.Range("D16").Value = 0
$salineSheet.Range("E16").Value = 0
$salineSheet.Range("F16").Value = "TF"

# Row 17: TF_good_placeholder
$salineSheet.Range("A17").Value = "TF_good_placeholder"
$salineSheet.Range("B17").Value = 0
$salineSheet.Range("C17").Value = 0.15
$salineSheet.Range("D17").Value = 1.93
$salineSheet.Range("E17").Value = 3
$salineSheet.Range("F17").Value = "TF"

# Column A autofits to the new, longer placeholder labels.
$salineSheet.Columns.Item(1).AutoFit() | Out-Null

# --- View state updates observed in the diff ---
# "Saline Sites" becomes the selected/active tab with a new selection.
# (Note: selecting a range on the other, non-active sheet would implicitly
# re-activate it, so we leave "Tidal Fresh Sites" untouched here -- it loses
# tabSelected purely as a side effect of "Saline Sites" becoming active.)
$salineSheet.Activate()
$salineSheet.Range("G18").Select() | Out-Null
